$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# Update the panel_query_time-equivalent cell (time_taken) on the "data" sheet
$dataSheet.Range("F2").Value = "2021-10-05 14:35:56.130728"

# Add a new "metadata" worksheet right after "data"
$newSheet = $wb.Worksheets.Add($null, $dataSheet)
$newSheet.Name = "metadata"

# Copy the header/row-2 cell formatting from the "data" sheet so the new
# sheet reuses the same existing styles (bold/bordered header, etc.)
$dataSheet.Range("B1").Copy()
$newSheet.Range("B1:G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row
$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

# Data row
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "Ventricular Fibrillation"
$newSheet.Range("C2").Value = 183
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "0.4"
$newSheet.Range("D2").Style = "Normal"
$newSheet.Range("E2").Value = "2021-04-06T10:32:33.119966Z"
$newSheet.Range("F2").Value = "2021-10-05 14:35:56.126978"
$newSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/183/?format=json"

# Match the page margins used by the "data" sheet
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Leave "data" as the active/selected sheet, as it was before the edit
$dataSheet.Activate()
